$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 16.54545838086715
$ws.Range("D2").Value = 5.722038895503942
$ws.Range("E2").Value = 17.52249396634937
$ws.Range("F2").Value = 29.5519018095289
$ws.Range("G2").Value = 38.08112812468947
$ws.Range("H2").Value = 16.42631398368536
$ws.Range("K2").Value = 9.318829144821382
$ws.Range("L2").Value = 8.75353456115476
$ws.Range("M2").Value = 15.43853456028563
$ws.Range("N2").Value = 21.62904568115881
$ws.Range("B3").Value = 16.47607657064935
$ws.Range("D3").Value = 5.711757341918568
$ws.Range("E3").Value = 17.5429456080166
$ws.Range("F3").Value = 29.48370385422345
$ws.Range("G3").Value = 37.93385562560665
$ws.Range("H3").Value = 16.45192886621513
$ws.Range("K3").Value = 8.893747222936744
$ws.Range("L3").Value = 8.740986830349208
$ws.Range("M3").Value = 15.42873960270313
$ws.Range("N3").Value = 21.69320454541203
$ws.Range("B4").Value = 16.4367342744358
$ws.Range("D4").Value = 5.7053036584773
$ws.Range("E4").Value = 17.55632372529484
$ws.Range("F4").Value = 29.44954628550543
$ws.Range("G4").Value = 37.85460718003018
$ws.Range("H4").Value = 16.47065086303179
$ws.Range("K4").Value = 8.620037964399353
$ws.Range("L4").Value = 8.734732667026895
$ws.Range("M4").Value = 15.4251837801056
$ws.Range("N4").Value = 21.73449619459756
$ws.Range("B5").Value = 16.42153428103252
$ws.Range("D5").Value = 5.702638675374615
$ws.Range("E5").Value = 17.56198224081656
$ws.Range("F5").Value = 29.43757592833875
$ws.Range("G5").Value = 37.82514614011323
$ws.Range("H5").Value = 16.47903228706476
$ws.Range("K5").Value = 8.505379384193407
$ws.Range("L5").Value = 8.732550804357409
$ws.Range("M5").Value = 15.42435535234353
$ws.Range("N5").Value = 21.75180131725315
$ws.Range("B6").Value = 16.41906096174855
$ws.Range("D6").Value = 5.702194049778452
$ws.Range("E6").Value = 17.562934337849
$ws.Range("F6").Value = 29.43570619898383
$ws.Range("G6").Value = 37.8204258659943
$ws.Range("H6").Value = 16.48046940947169
$ws.Range("K6").Value = 8.486154442707718
$ws.Range("L6").Value = 8.732210715782985
$ws.Range("M6").Value = 15.42425534053126
$ws.Range("N6").Value = 21.7547037520024
$ws.Range("B7").Value = 16.43652589566106
$ws.Range("D7").Value = 5.705267859074935
$ws.Range("E7").Value = 17.55639919997604
$ws.Range("F7").Value = 29.44937694708531
$ws.Range("G7").Value = 37.8541983591988
$ws.Range("H7").Value = 16.47076085411566
$ws.Range("K7").Value = 8.618504161657556
$ws.Range("L7").Value = 8.734701754074317
$ws.Range("M7").Value = 15.4251700918654
$ws.Range("N7").Value = 21.73472763875407
$ws.Range("B8").Value = 16.52086732681999
$ws.Range("D8").Value = 5.71852323547941
$ws.Range("E8").Value = 17.52937561834223
$ws.Range("F8").Value = 29.52679155523882
$ws.Range("G8").Value = 38.02804352928803
$ws.Range("H8").Value = 16.4345240189729
$ws.Range("K8").Value = 9.174939748218105
$ws.Range("L8").Value = 8.748908440764815
$ws.Range("M8").Value = 15.43464857818663
$ws.Range("N8").Value = 21.65077449162876
$ws.Range("B9").Value = 16.71149990497549
$ws.Range("D9").Value = 5.743385807945186
$ws.Range("E9").Value = 17.48287455432732
$ws.Range("F9").Value = 29.73934640557311
$ws.Range("G9").Value = 38.45645484262636
$ws.Range("H9").Value = 16.38725798303938
$ws.Range("K9").Value = 10.16256797246201
$ws.Range("L9").Value = 8.788173131899402
$ws.Range("M9").Value = 15.47262384694731
$ws.Range("N9").Value = 21.50114488900553
$ws.Range("B10").Value = 16.86605825209355
$ws.Range("D10").Value = 5.760944785562759
$ws.Range("E10").Value = 17.45264047579068
$ws.Range("F10").Value = 29.93176654692595
$ws.Range("G10").Value = 38.82269668722251
$ws.Range("H10").Value = 16.36707866278614
$ws.Range("K10").Value = 10.82214694208658
$ws.Range("L10").Value = 8.823825098478357
$ws.Range("M10").Value = 15.51217012768665
$ws.Range("N10").Value = 21.40027754391105
$ws.Range("B11").Value = 16.93931270893779
$ws.Range("D11").Value = 5.76877543445878
$ws.Range("E11").Value = 17.43973385833387
$ws.Range("F11").Value = 30.02696865203248
$ws.Range("G11").Value = 39.00001631292432
$ws.Range("H11").Value = 16.3610628816449
$ws.Range("K11").Value = 11.10742570140577
$ws.Range("L11").Value = 8.841485118437657
$ws.Range("M11").Value = 15.53264642466565
$ws.Range("N11").Value = 21.35634181300065
$ws.Range("B12").Value = 16.96745714170027
$ws.Range("D12").Value = 5.771717776223228
$ws.Range("E12").Value = 17.43496782286604
$ws.Range("F12").Value = 30.06410142049801
$ws.Range("G12").Value = 39.06865820920918
$ws.Range("H12").Value = 16.35923994041209
$ws.Range("K12").Value = 11.21330169089116
$ws.Range("L12").Value = 8.848376358528489
$ws.Range("M12").Value = 15.54075363343172
$ws.Range("N12").Value = 21.33998356806795
$ws.Range("B13").Value = 16.96137805023913
$ws.Range("D13").Value = 5.771085117180906
$ws.Range("E13").Value = 17.43598887965972
$ws.Range("F13").Value = 30.05605645560297
$ws.Range("G13").Value = 39.05380923090115
$ws.Range("H13").Value = 16.35961230229381
$ws.Range("K13").Value = 11.19059555952862
$ws.Range("L13").Value = 8.846883202637505
$ws.Range("M13").Value = 15.53899196003881
$ws.Range("N13").Value = 21.34349420504473
$ws.Range("B14").Value = 16.94162017122573
$ws.Range("D14").Value = 5.769017964889601
$ws.Range("E14").Value = 17.43933932243251
$ws.Range("F14").Value = 30.03000203342071
$ws.Range("G14").Value = 39.0056338389975
$ws.Range("H14").Value = 16.36090378606772
$ws.Range("K14").Value = 11.11617947774455
$ws.Range("L14").Value = 8.84204800209069
$ws.Range("M14").Value = 15.53330635492029
$ws.Range("N14").Value = 21.35499041945584
$ws.Range("B15").Value = 16.92957001194291
$ws.Range("D15").Value = 5.767748773860985
$ws.Range("E15").Value = 17.44140736673828
$ws.Range("F15").Value = 30.01418317381248
$ws.Range("G15").Value = 38.97631827197516
$ws.Range("H15").Value = 16.3617541261602
$ws.Range("K15").Value = 11.07031633599046
$ws.Range("L15").Value = 8.839112730396383
$ws.Range("M15").Value = 15.52986963301039
$ws.Range("N15").Value = 21.36206852012561
$ws.Range("B16").Value = 16.86132856022508
$ws.Range("D16").Value = 5.760429842465223
$ws.Range("E16").Value = 17.45350096516351
$ws.Range("F16").Value = 29.92569743494184
$ws.Range("G16").Value = 38.81132018077609
$ws.Range("H16").Value = 16.36753549015927
$ws.Range("K16").Value = 10.80320350548761
$ws.Range("L16").Value = 8.822699688186518
$ws.Range("M16").Value = 15.51088162957828
$ws.Range("N16").Value = 21.40318798279489
$ws.Range("B17").Value = 16.82020569083995
$ws.Range("D17").Value = 5.755899431120513
$ws.Range("E17").Value = 17.46113666547318
$ws.Range("F17").Value = 29.87336420686633
$ws.Range("G17").Value = 38.71281175221737
$ws.Range("H17").Value = 16.37189267827286
$ws.Range("K17").Value = 10.63553342057365
$ws.Range("L17").Value = 8.812997633143247
$ws.Range("M17").Value = 15.49986700969423
$ws.Range("N17").Value = 21.42891196084791
$ws.Range("B18").Value = 16.79683128915014
$ws.Range("D18").Value = 5.753278901434637
$ws.Range("E18").Value = 17.4656082659308
$ws.Range("F18").Value = 29.84398676132446
$ws.Range("G18").Value = 38.65716375497203
$ws.Range("H18").Value = 16.37469662152933
$ws.Range("K18").Value = 10.53770543005985
$ws.Range("L18").Value = 8.807553199827415
$ws.Range("M18").Value = 15.4937660075379
$ws.Range("N18").Value = 21.44389125452065
$ws.Range("B19").Value = 16.78896549717127
$ws.Range("D19").Value = 5.752389112082742
$ws.Range("E19").Value = 17.46713598137776
$ws.Range("F19").Value = 29.8341648907721
$ws.Range("G19").Value = 38.63849736289644
$ws.Range("H19").Value = 16.3756971293137
$ws.Range("K19").Value = 10.50434505477183
$ws.Range("L19").Value = 8.805733259254202
$ws.Range("M19").Value = 15.49174068234339
$ws.Range("N19").Value = 21.44899454168684
$ws.Range("B20").Value = 16.82455461372133
$ws.Range("D20").Value = 5.756383231144458
$ws.Range("E20").Value = 17.46031558134847
$ws.Range("F20").Value = 29.87886046174146
$ws.Range("G20").Value = 38.72319378431393
$ws.Range("H20").Value = 16.37139802493138
$ws.Range("K20").Value = 10.65352613301406
$ws.Range("L20").Value = 8.814016390008984
$ws.Range("M20").Value = 15.50101531398083
$ws.Range("N20").Value = 21.42615461119603
$ws.Range("B21").Value = 16.94741271382851
$ws.Range("D21").Value = 5.769625763609247
$ws.Range("E21").Value = 17.43835192385502
$ws.Range("F21").Value = 30.0376256658459
$ws.Range("G21").Value = 39.01974393217083
$ws.Range("H21").Value = 16.36051209438614
$ws.Range("K21").Value = 11.13809591927744
$ws.Range("L21").Value = 8.843462716715344
$ws.Range("M21").Value = 15.53496680262387
$ws.Range("N21").Value = 21.35160613177928
$ws.Range("B22").Value = 17.03005643253165
$ws.Range("D22").Value = 5.778146551093116
$ws.Range("E22").Value = 17.42470496330355
$ws.Range("F22").Value = 30.14768241147894
$ws.Range("G22").Value = 39.22224602159711
$ws.Range("H22").Value = 16.35605013805585
$ws.Range("K22").Value = 11.44223188503441
$ws.Range("L22").Value = 8.863893543155662
$ws.Range("M22").Value = 15.55921305244971
$ws.Range("N22").Value = 21.30451167249169
$ws.Range("B23").Value = 16.98573942938782
$ws.Range("D23").Value = 5.773611228910145
$ws.Range("E23").Value = 17.43192398512654
$ws.Range("F23").Value = 30.08837452650042
$ws.Range("G23").Value = 39.11338771889457
$ws.Range("H23").Value = 16.35818885500905
$ws.Range("K23").Value = 11.28106615144421
$ws.Range("L23").Value = 8.852881941204275
$ws.Range("M23").Value = 15.54608564772862
$ws.Range("N23").Value = 21.32949831796986
$ws.Range("B24").Value = 16.82258763044072
$ws.Range("D24").Value = 5.756164554745824
$ws.Range("E24").Value = 17.4606865388912
$ws.Range("F24").Value = 29.87637339284781
$ws.Range("G24").Value = 38.71849699334675
$ws.Range("H24").Value = 16.37162072650537
$ws.Range("K24").Value = 10.64539608658111
$ws.Range("L24").Value = 8.813555394210239
$ws.Range("M24").Value = 15.50049544434332
$ws.Range("N24").Value = 21.42740061628436
$ws.Range("B25").Value = 16.65731457256375
$ws.Range("D25").Value = 5.736782780777101
$ws.Range("E25").Value = 17.49476216441191
$ws.Range("F25").Value = 29.67541344218705
$ws.Range("G25").Value = 38.3313669364329
$ws.Range("H25").Value = 16.39749241123122
$ws.Range("K25").Value = 9.906808795603979
$ws.Range("L25").Value = 8.776344393853845
$ws.Range("M25").Value = 15.46029161342575
$ws.Range("N25").Value = 21.54002591493527
